$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: experience 0.9 -> 0.85, final_score 75.81999999999999 -> 71.59999999999999
$ws.Range("E3").Value = 0.85
$ws.Range("F3").Value = 71.59999999999999

# Row 6: experience 0.9 -> 0.85, final_score 66.83 -> 63.11
$ws.Range("E6").Value = 0.85
$ws.Range("F6").Value = 63.11

# Row 8: experience 0.9 -> 0.85, final_score 72.77 -> 68.72
$ws.Range("E8").Value = 0.85
$ws.Range("F8").Value = 68.72
